# The commit swaps the presentation's theme palette: the deck's design
# ("Integral" / "Red Violet") is replaced by the default "Office Theme"
# color scheme (the palette that, before this edit, only the notes
# master theme part carried).
#
# Helper: turn an "RRGGBB" hex string into the OLE RGB() long that the
# PowerPoint object model's ColorFormat/ThemeColor .RGB property expects
# (r + g*256 + b*65536).
function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgb $officeTheme[$i - 1]
}
